# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet is missing a "property_category" column
# that every other asset-type sheet (land, building, car, ...) implicitly
# carries via its sheet name. Insert that column right after "total" and
# before "date" (i.e. at column H), and stamp every data row with the
# literal value "stock".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Shift the existing date / legislator_name / legislator_id columns one
# slot to the right, opening up column H for the new field.
$ws.Columns("H:H").Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H10").Value = "stock"
